$d = $word.ActiveDocument

# =====================================================================
# Change 1: remove the stray "_GoBack" bookmark that currently sits in
# the "Add the following items ... to Java Options for Tomcat" paragraph
# =====================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# =====================================================================
# Change 2: collapse the three runs that make up the "* Note that you
# can use a directory other than C:\shindig\conf, ..." sentence into a
# single run (this also fixes the word "shindig" which was split across
# two runs as "shind" + "ig").
# =====================================================================
$oldNote = "* Note that you can use a directory other than C:\shind" + "ig\conf, but whatever you use you must make sure it is in the beginning of the Java "
$newNote = "* Note that you can use a directory other than C:\shindig\conf, but whatever you use you must make sure it is in the beginning of the Java "

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($oldNote, $true, $true, $false, $false, $false, $true, 1, $false, $newNote, 2) | Out-Null

# =====================================================================
# Change 3: add a new "9. Hosting Gadgets" Heading-3 section (with a
# paragraph of body text) right before the existing "9. Join Us and
# Contribute" heading, and move the "_GoBack" bookmark to the end of
# the new body paragraph.
# =====================================================================
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "9. Join Us*") {
        $target = $p
        break
    }
}

$insertionPoint = $target.Range
$insertionPoint.Collapse(1) | Out-Null   # wdCollapseStart
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()

# Re-fetch the two freshly-created (still empty) paragraphs.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "9. Join Us*") {
        $target = $p
        break
    }
}
$headingPara = $target.Previous(2)
$headingPara2 = $headingPara.Next(1)

$headingPara.Style = "Heading 3"
$headingPara2.Style = "Normal"

$headingRange = $headingPara.Range
$headingRange.MoveEnd(1, -1) | Out-Null
$headingRange.Text = "9. Hosting Gadgets"

$bodyText = "You will notice that the default gadgets are hosted on external web sites.  Feel free to use these where they are, or to copy them to one of your own web servers where you can modify them if desired. At UCSF we host them on the same IIS web servers that we use for Profiles."

$bodyRange = $headingPara2.Range
$bodyRange.MoveEnd(1, -1) | Out-Null
$bodyRange.Text = $bodyText

# ---- place the "_GoBack" bookmark right after the body text, at the
# ---- very end of the paragraph (before the paragraph mark). Adding a
# ---- bookmark exactly at the paragraph-end boundary is unreliable, so
# ---- append a throw-away character, anchor the bookmark just before
# ---- it, then remove the throw-away character again.
$bodyPara = $headingPara2
$tempRange = $bodyPara.Range
$tempRange.MoveEnd(1, -1) | Out-Null
$tempRange.InsertAfter("X")

$bodyPara = $headingPara2
$fullRange = $bodyPara.Range
$fullRange.MoveEnd(1, -1) | Out-Null
$bmPos = $fullRange.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$bodyPara = $headingPara2
$cleanupRange = $bodyPara.Range
$cleanupRange.MoveEnd(1, -1) | Out-Null
$lastCharRange = $cleanupRange.Duplicate()
$lastCharRange.MoveStart(1, $lastCharRange.End - $lastCharRange.Start - 1) | Out-Null
$lastCharRange.Delete()
